# update content according to prof. comment
#
# This script applies two related edits to the "research projects"
# paragraph and relocates the lone "_GoBack" bookmark that Word leaves
# behind at the last edit point:
#
#  1. In the bullet paragraph that reads
#       "Participate in two research project: a MOST project on
#        multimedia and an UMC joint project on machine learning"
#     - merge the split "...research " / "project" runs (and drop the
#       now-redundant proofErr gramStart/gramEnd markers around "project"),
#     - pluralize "project" -> "projects",
#     - add a trailing ".",
#     - and plant the "_GoBack" bookmark right after "projects".
#
#  2. Remove the old "_GoBack" bookmark that used to sit after
#     "本獎學金申請表" (it has moved to the paragraph above).
#
# Because the COM Font object in this host only round-trips
# ascii/hAnsi font names (eastAsia/cs/hint are not settable through
# Font.Name / Font.NameFarEast / Font.NameOther), the precise run
# splitting + w:hint="eastAsia" marking required by the target markup
# is produced with Range.InsertXML, rebuilding just the two affected
# paragraphs from a FlatOpc fragment while leaving every other part of
# the document untouched.

$d = $word.ActiveDocument

function Get-ContainingParagraphRange($doc, [int]$pos) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $pr = $p.Range
        if ($pr.Start -le $pos -and $pr.End -gt $pos) {
            return $pr
        }
    }
    throw "no paragraph contains position $pos"
}

function Get-ParagraphTextRange($doc, $paraRange) {
    # Paragraph.Range includes the trailing paragraph mark (chr 13) and,
    # for the last paragraph in a table cell, the cell mark (chr 7) as
    # well. Trim those off so InsertXML only replaces the paragraph's
    # own content (and keeps it as a single paragraph).
    $endPos = $paraRange.End
    while ($endPos -gt $paraRange.Start) {
        $lastChar = $doc.Range($endPos - 1, $endPos).Text
        $code = [int][char]$lastChar
        if ($code -eq 13 -or $code -eq 7) {
            $endPos = $endPos - 1
        } else {
            break
        }
    }
    return $doc.Range($paraRange.Start, $endPos)
}

function New-FlatOpcPayload([string]$paragraphXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $paragraphXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# Step 1 — drop the stray "_GoBack" bookmark after "本獎學金申請表".
# Do this first: it does not change any character offsets (a bookmark
# is zero-width), but rebuilding it via InsertXML before step 2 keeps
# every later offset we still care about (the earlier "research
# project..." paragraph) untouched regardless of ordering.
# ---------------------------------------------------------------------

$find1 = $d.Content
$found1 = $find1.Find.Execute("本獎學金申請表", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "could not find the '本獎學金申請表' paragraph" }

$para1Range = Get-ContainingParagraphRange $d $find1.Start
$target1 = Get-ParagraphTextRange $d $para1Range

$checkbox = [char]9633
$para1Xml = '<w:p w:rsidR="00431487" w:rsidRPr="00431487" w:rsidRDefault="00431487" w:rsidP="00431487">' +
    '<w:pPr><w:ind w:left="357" w:hanging="357"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:cs="Times New Roman"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="007B7C90"><w:rPr><w:rFonts w:ascii="新細明體" w:eastAsia="標楷體" w:hAnsi="新細明體" w:cs="Times New Roman" w:hint="eastAsia"/><w:highlight w:val="darkGray"/></w:rPr><w:t>' + $checkbox + '</w:t></w:r>' +
    '<w:r w:rsidRPr="00431487"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>本獎學金申請表</w:t></w:r>' +
    '</w:p>'

$target1.InsertXML((New-FlatOpcPayload $para1Xml))

# ---------------------------------------------------------------------
# Step 2 — rebuild the "Participate in two research project(s)" bullet:
# merge the runs, pluralize, add the trailing period, and re-home the
# "_GoBack" bookmark right after "projects".
# ---------------------------------------------------------------------

$find2 = $d.Content
$found2 = $find2.Find.Execute("Participate in two research project", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "could not find the 'Participate in two research project' paragraph" }

$para2Range = Get-ContainingParagraphRange $d $find2.Start
$target2 = Get-ParagraphTextRange $d $para2Range

$para2Xml = '<w:p w:rsidR="007B7C90" w:rsidRPr="007B7C90" w:rsidRDefault="007B7C90" w:rsidP="007B7C90">' +
    '<w:pPr><w:pStyle w:val="ac"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:cs="Times New Roman"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>P</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:cs="Times New Roman"/></w:rPr><w:t>articipate in two research project</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:cs="Times New Roman"/></w:rPr><w:t>s</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:cs="Times New Roman"/></w:rPr><w:t>: a MOST project on multimedia and an UMC joint project on machine learning</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>.</w:t></w:r>' +
    '</w:p>'

$target2.InsertXML((New-FlatOpcPayload $para2Xml))

Write-Output "edit applied"
